# Auto-generated edit script: applies cell-value corrections from the
# authoritative diff to the Exodus_Profits workbook (8 profession sheets).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 173.5
$ws.Range("I33").Value = 199.3077
$ws.Range("J33").Value = 61.666668
$ws.Range("K33").Value = 199.3077
$ws.Range("L33").Value = 61.666668
$ws.Range("M33").Value = 29.69229999999999
$ws.Range("N33").Value = -519.666668
$ws.Range("H74").Value = 5756
$ws.Range("I74").Value = 5779.8
$ws.Range("K74").Value = 5779.8
$ws.Range("M74").Value = -4843.8
$ws.Range("H77").Value = 5756
$ws.Range("I77").Value = 5779.8
$ws.Range("K77").Value = 28899
$ws.Range("M77").Value = -24219
$ws.Range("H132").Value = 1607
$ws.Range("I132").Value = 1679.2122
$ws.Range("J132").Value = 1011.25
$ws.Range("K132").Value = 5037.6366
$ws.Range("L132").Value = 3033.75
$ws.Range("M132").Value = -2507.6366
$ws.Range("N132").Value = -8093.75
$ws.Range("H138").Value = 16511.492
$ws.Range("I138").Value = 56891.777
$ws.Range("J138").Value = 2259.6274
$ws.Range("K138").Value = 170675.331
$ws.Range("L138").Value = 6778.8822
$ws.Range("M138").Value = -165535.331
$ws.Range("N138").Value = -17058.8822

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 1422.25
$ws.Range("I14").Value = 250
$ws.Range("J14").Value = 2594.5
$ws.Range("K14").Value = 250
$ws.Range("L14").Value = 2594.5
$ws.Range("M14").Value = -75
$ws.Range("N14").Value = -2944.5
$ws.Range("H27").Value = 2166.6667
$ws.Range("J27").Value = 2166.6667
$ws.Range("L27").Value = 2166.6667
$ws.Range("N27").Value = -2534.6667
$ws.Range("H32").Value = 7055.3276
$ws.Range("I32").Value = 3884.1428
$ws.Range("K32").Value = 3884.1428
$ws.Range("M32").Value = -3597.1428
$ws.Range("H61").Value = 35626.3
$ws.Range("I61").Value = 2065.1304
$ws.Range("K61").Value = 2065.1304
$ws.Range("M61").Value = -1853.1304
$ws.Range("H74").Value = 42427.12
$ws.Range("I74").Value = 63964.562
$ws.Range("K74").Value = 63964.562
$ws.Range("M74").Value = -63090.562
$ws.Range("H77").Value = 42427.12
$ws.Range("I77").Value = 63964.562
$ws.Range("K77").Value = 319822.81
$ws.Range("M77").Value = -315454.81
$ws.Range("H136").Value = 35626.3
$ws.Range("I136").Value = 2065.1304
$ws.Range("K136").Value = 6195.3912
$ws.Range("M136").Value = -3645.3912

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4112.136
$ws.Range("I86").Value = 2972.5625
$ws.Range("K86").Value = 2972.5625
$ws.Range("M86").Value = -1849.5625
$ws.Range("H89").Value = 4112.136
$ws.Range("I89").Value = 2972.5625
$ws.Range("K89").Value = 14862.8125
$ws.Range("M89").Value = -9246.8125
$ws.Range("H94").Value = 1642.8636
$ws.Range("I94").Value = 1407.4445
$ws.Range("K94").Value = 1407.4445
$ws.Range("M94").Value = -956.4445000000001
$ws.Range("H105").Value = 48745.184
$ws.Range("I105").Value = 127449.375
$ws.Range("J105").Value = 3771.3572
$ws.Range("K105").Value = 127449.375
$ws.Range("L105").Value = 3771.3572
$ws.Range("M105").Value = -125702.375
$ws.Range("N105").Value = -7265.3572
$ws.Range("H134").Value = 2663.9524
$ws.Range("I134").Value = 1068.3125
$ws.Range("J134").Value = 7770
$ws.Range("K134").Value = 3204.9375
$ws.Range("L134").Value = 23310
$ws.Range("M134").Value = -669.9375
$ws.Range("N134").Value = -28380

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 619.6667
$ws.Range("I10").Value = 225.28572
$ws.Range("K10").Value = 225.28572
$ws.Range("M10").Value = -86.28572
$ws.Range("H31").Value = 2177.9487
$ws.Range("I31").Value = 1635.5807
$ws.Range("J31").Value = 4279.625
$ws.Range("K31").Value = 1635.5807
$ws.Range("L31").Value = 4279.625
$ws.Range("M31").Value = -1340.5807
$ws.Range("N31").Value = -4869.625
$ws.Range("H34").Value = 2177.9487
$ws.Range("I34").Value = 1635.5807
$ws.Range("J34").Value = 4279.625
$ws.Range("K34").Value = 1635.5807
$ws.Range("L34").Value = 4279.625
$ws.Range("M34").Value = -1433.5807
$ws.Range("N34").Value = -4683.625
$ws.Range("H99").Value = 4833269
$ws.Range("I99").Value = 6946622.5
$ws.Range("J99").Value = 2746.7144
$ws.Range("K99").Value = 6946622.5
$ws.Range("L99").Value = 2746.7144
$ws.Range("M99").Value = -6945124.5
$ws.Range("N99").Value = -5742.7144
$ws.Range("H105").Value = 4609.8
$ws.Range("I105").Value = 250
$ws.Range("K105").Value = 250
$ws.Range("M105").Value = 1497
$ws.Range("H107").Value = 698.0833
$ws.Range("I107").Value = 696.4
$ws.Range("K107").Value = 696.4
$ws.Range("M107").Value = 1223.6
$ws.Range("H126").Value = 4833269
$ws.Range("I126").Value = 6946622.5
$ws.Range("J126").Value = 2746.7144
$ws.Range("K126").Value = 20839867.5
$ws.Range("L126").Value = 8240.143199999999
$ws.Range("M126").Value = -20837397.5
$ws.Range("N126").Value = -13180.1432

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 102.875
$ws.Range("I2").Value = 181.5
$ws.Range("K2").Value = 1089
$ws.Range("M2").Value = -976
$ws.Range("H38").Value = 489.81818
$ws.Range("J38").Value = 743.1429000000001
$ws.Range("L38").Value = 2229.4287
$ws.Range("N38").Value = -2923.4287
$ws.Range("H56").Value = 6697
$ws.Range("I56").Value = 6697
$ws.Range("K56").Value = 6697
$ws.Range("M56").Value = -6167
$ws.Range("H131").Value = 1240.1538
$ws.Range("I131").Value = 1032
$ws.Range("K131").Value = 3096
$ws.Range("M131").Value = 1944
$ws.Range("H141").Value = 1648.5454
$ws.Range("I141").Value = 1648.5454
$ws.Range("K141").Value = 4945.6362
$ws.Range("M141").Value = 234.3638000000001

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 112137.38
$ws.Range("I70").Value = 60197.445
$ws.Range("K70").Value = 60197.445
$ws.Range("M70").Value = -59927.445
$ws.Range("H73").Value = 112137.38
$ws.Range("I73").Value = 60197.445
$ws.Range("K73").Value = 60197.445
$ws.Range("M73").Value = -59261.445
$ws.Range("H132").Value = 4919
$ws.Range("I132").Value = 2768.3333
$ws.Range("K132").Value = 8304.999899999999
$ws.Range("M132").Value = -5774.999899999999
$ws.Range("H135").Value = 53844.75
$ws.Range("J135").Value = 53844.75
$ws.Range("L135").Value = 53844.75
$ws.Range("N135").Value = -63984.75

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3655.6858
$ws.Range("I22").Value = 688.64703
$ws.Range("J22").Value = 6457.8887
$ws.Range("K22").Value = 688.64703
$ws.Range("L22").Value = 6457.8887
$ws.Range("M22").Value = -393.64703
$ws.Range("N22").Value = -7047.8887
$ws.Range("H27").Value = 3655.6858
$ws.Range("I27").Value = 688.64703
$ws.Range("J27").Value = 6457.8887
$ws.Range("K27").Value = 688.64703
$ws.Range("L27").Value = 6457.8887
$ws.Range("M27").Value = -581.64703
$ws.Range("N27").Value = -6671.8887
$ws.Range("H40").Value = 3707408
$ws.Range("I40").Value = 3825
$ws.Range("K40").Value = 3825
$ws.Range("M40").Value = -3689
$ws.Range("H82").Value = 4497.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 4497.5
$ws.Range("K82").Value = 0
$ws.Range("L82").ClearContents() | Out-Null
$ws.Range("M82").Value = 4497.5
$ws.Range("N82").Value = -5219.5
$ws.Range("H85").Value = 4497.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 4497.5
$ws.Range("K85").Value = 0
$ws.Range("L85").ClearContents() | Out-Null
$ws.Range("M85").Value = 4497.5
$ws.Range("N85").Value = -6993.5
$ws.Range("H122").Value = 92312024
$ws.Range("I122").Value = 90913660
$ws.Range("J122").Value = 100003030
$ws.Range("K122").Value = 272740980
$ws.Range("L122").Value = 300009090
$ws.Range("M122").Value = -272738530
$ws.Range("N122").Value = -300013990
$ws.Range("H134").Value = 127862
$ws.Range("J134").Value = 127862
$ws.Range("L134").Value = 127862
$ws.Range("N134").Value = -138002
$ws.Range("H136").Value = 3933.92
$ws.Range("I136").Value = 4972.091
$ws.Range("J136").Value = 3118.2144
$ws.Range("K136").Value = 14916.273
$ws.Range("L136").Value = 9354.643199999999
$ws.Range("M136").Value = -12366.273
$ws.Range("N136").Value = -14454.6432

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6210200
$ws.Range("I2").Value = 10016666
$ws.Range("K2").Value = 10016666
$ws.Range("M2").Value = -10016554
$ws.Range("H4").Value = 7329.143
$ws.Range("I4").Value = 2150
$ws.Range("K4").Value = 2150
$ws.Range("M4").Value = -2037
$ws.Range("H59").Value = 33499.375
$ws.Range("J59").Value = 33499.375
$ws.Range("L59").Value = 33499.375
$ws.Range("N59").Value = -34975.375
$ws.Range("H61").Value = 7082171
$ws.Range("I61").Value = 9253527
$ws.Range("J61").Value = 25263
$ws.Range("K61").Value = 9253527
$ws.Range("L61").Value = 25263
$ws.Range("M61").Value = -9253235
$ws.Range("N61").Value = -25847
$ws.Range("H122").Value = 3912.8
$ws.Range("I122").Value = 2036.5
$ws.Range("K122").Value = 6109.5
$ws.Range("M122").Value = -3659.5
$ws.Range("H126").Value = 4384.3335
$ws.Range("I126").Value = 3269
$ws.Range("K126").Value = 9807
$ws.Range("M126").Value = -7337
$ws.Range("H132").Value = 2073.25
$ws.Range("I132").Value = 1710.909
$ws.Range("K132").Value = 5132.727000000001
$ws.Range("M132").Value = -2602.727000000001
$ws.Range("H136").Value = 1706.0454
$ws.Range("I136").Value = 1457.1666
$ws.Range("J136").Value = 2004.7
$ws.Range("K136").Value = 4371.4998
$ws.Range("L136").Value = 6014.1
$ws.Range("M136").Value = -1821.4998
$ws.Range("N136").Value = -11114.1

"Applied 248 cell updates and 2 clears."
